$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P4: new year header (2022), formatted like the other year headers (O4) ---
$ws.Range("P4").Value = 2022
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)

# --- P5: total row, bold font matching column A "Total" style ---
$r5 = $ws.Range("P5")
$r5.Value = 1188.7
$r5.HorizontalAlignment = -4152
$r5.NumberFormat = "#,##0.0"
$r5.VerticalAlignment = -4107
$r5.Font.Bold = $true
$r5.Font.Color = 0

# --- P6:P15: regular data rows, right aligned, new thousands-separated format ---
$vals = @{6=263.9; 7=263.2; 8=12.4; 9="-"; 10=93; 11=171.5; 12=220.6; 13=159.3; 14=1.7; 15="-"}
foreach ($row in 6..15) {
    $ws.Cells.Item($row, 16).Value = $vals[$row]
}
$rng = $ws.Range("P6:P15")
$rng.HorizontalAlignment = -4152
$rng.NumberFormat = "#,##0.0"
$rng.VerticalAlignment = -4107

# --- P16: bottom (footer) row, keep the medium bottom border from O16 ---
$ws.Range("O16").Copy()
$ws.Range("P16").PasteSpecial(-4122)
$r16 = $ws.Range("P16")
$r16.Value = 3.1
$r16.NumberFormat = "#,##0.0"
$r16.VerticalAlignment = -4107

# --- Selection, to match the cursor position left behind by the edit ---
$ws.Range("Q7").Select() | Out-Null
